$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item('ランサーズ')

# Clear any existing hyperlinks on the sheet so stale relationship targets
# are not left pointing at the old URLs once cell text is rewritten.
$ws.Hyperlinks.Delete()

# Row 2
$ws.Range('A2').Value = '2026-02-16 18:39:57'
$ws.Range('B2').Value = '製造業向け図面自動生成システムの開発・ツール化を支援してくださるエンジニア募集(AI/バックエンド)'
$ws.Range('C2').Value = 'システム開発'
$ws.Range('D2').Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range('E2').Value = '期限情報なし'
$ws.Range('F2').Value = 'https://www.lancers.jp/work/detail/5460562'
$ws.Range('G2').Value = 435
$ws.Range('H2').Value = '🔥AI,Ai ◆ツール,開発'

# Row 3
$ws.Range('A3').Value = '2026-02-16 18:39:57'
$ws.Range('B3').Value = '製造業向け設備要件定義書の自動生成AIシステムの開発・DB設計支援エンジニア(AI/バックエンド)'
$ws.Range('C3').Value = 'システム開発'
$ws.Range('D3').Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range('E3').Value = '期限情報なし'
$ws.Range('F3').Value = 'https://www.lancers.jp/work/detail/5473648'
$ws.Range('G3').Value = 390
$ws.Range('H3').Value = '🔥AI,Ai ◆開発'

# Row 4
$ws.Range('A4').Value = '2026-02-16 18:39:57'
$ws.Range('B4').Value = '【急募】AI技術に精通したフリーランスを探しています!'
$ws.Range('C4').Value = 'システム開発'
$ws.Range('D4').Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range('E4').Value = '期限情報なし'
$ws.Range('F4').Value = 'https://www.lancers.jp/work/detail/5492832'
$ws.Range('G4').Value = 310
$ws.Range('H4').Value = '🔥AI,Ai'

# Row 5
$ws.Range('A5').Value = '2026-02-16 18:39:57'
$ws.Range('B5').Value = '【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪'
$ws.Range('C5').Value = 'システム開発'
$ws.Range('D5').Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range('E5').Value = '期限情報なし'
$ws.Range('F5').Value = 'https://www.lancers.jp/work/detail/5217096'
$ws.Range('G5').Value = 243
$ws.Range('H5').Value = '🔥API ◆ツール'

# Row 6
$ws.Range('A6').Value = '2026-02-16 18:39:57'
$ws.Range('B6').Value = '施設管理・現場業務向け チェックリスト業務の自動化・報告書作成システム開発エンジニア募集'
$ws.Range('C6').Value = 'システム開発'
$ws.Range('D6').Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range('E6').Value = '期限情報なし'
$ws.Range('F6').Value = 'https://www.lancers.jp/work/detail/5460563'
$ws.Range('G6').Value = 220
$ws.Range('H6').Value = '◆開発,システム開発 ◇管理'

# Row 7
$ws.Range('A7').Value = '2026-02-16 18:39:57'
$ws.Range('B7').Value = '【Zapier保守・運用サポート】既存フローの管理・調整をお任せできる方募集(時給1,200円程度)'
$ws.Range('C7').Value = 'システム開発'
$ws.Range('D7').Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range('E7').Value = '期限情報なし'
$ws.Range('F7').Value = 'https://www.lancers.jp/work/detail/5488168'
$ws.Range('G7').Value = 213
$ws.Range('H7').Value = '🔥API ◇管理'

# Row 8
$ws.Range('A8').Value = '2026-02-16 18:39:57'
$ws.Range('B8').Value = '移行準備のためのファイルURL一覧(台帳)作成:API連携スクリプト+Excel作成'
$ws.Range('C8').Value = 'システム開発'
$ws.Range('D8').Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range('E8').Value = '期限情報なし'
$ws.Range('F8').Value = 'https://www.lancers.jp/work/detail/5492887'
$ws.Range('G8').Value = 188
$ws.Range('H8').Value = '🔥API'

# Row 9
$ws.Range('A9').Value = '2026-02-16 18:39:57'
$ws.Range('B9').Value = '【社内用】Temu APIを使って受注データ、在庫データ、注文ステータスを更新してほしい'
$ws.Range('C9').Value = 'システム開発'
$ws.Range('D9').Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range('E9').Value = '期限情報なし'
$ws.Range('F9').Value = 'https://www.lancers.jp/work/detail/5492576'
$ws.Range('G9').Value = 188
$ws.Range('H9').Value = '🔥API'

# Row 10
$ws.Range('A10').Value = '2026-02-16 18:39:57'
$ws.Range('B10').Value = '【買い切り20万円】Shopeeチャット管理・返信Webツール開発(複数国対応)'
$ws.Range('C10').Value = 'システム開発'
$ws.Range('D10').Value = '200,000 円 ~ 300,000 円 / 募集期間 3 日、取引期間 0 日'
$ws.Range('E10').Value = '期限情報なし'
$ws.Range('F10').Value = 'https://www.lancers.jp/work/detail/5492959'
$ws.Range('G10').Value = 163
$ws.Range('H10').Value = '◆ツール,開発 ◇管理'

# Row 11
$ws.Range('A11').Value = '2026-02-16 18:39:57'
$ws.Range('B11').Value = '【エンジニア募集】香水自販機制御システム開発'
$ws.Range('C11').Value = 'システム開発'
$ws.Range('D11').Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range('E11').Value = '期限情報なし'
$ws.Range('F11').Value = 'https://www.lancers.jp/work/detail/5492441'
$ws.Range('G11').Value = 125
$ws.Range('H11').Value = '◆開発,システム開発'

# Row 12
$ws.Range('A12').Value = '2026-02-16 18:39:57'
$ws.Range('B12').Value = '地域情報サイト 店舗データ自動収集・一括管理システム構築'
$ws.Range('C12').Value = 'システム開発'
$ws.Range('D12').Value = '1,000,000 円 ~ 3,000,000 円 / 固定'
$ws.Range('E12').Value = '期限情報なし'
$ws.Range('F12').Value = 'https://www.lancers.jp/work/detail/5492383'
$ws.Range('G12').Value = 85
$ws.Range('H12').Value = '◇サイト'

# Row 13
$ws.Range('A13').Value = '2026-02-16 18:39:57'
$ws.Range('B13').Value = '製造業DXプロダクト開発のプロダクトマネージャー募集'
$ws.Range('C13').Value = 'システム開発'
$ws.Range('D13').Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range('E13').Value = '期限情報なし'
$ws.Range('F13').Value = 'https://www.lancers.jp/work/detail/5468432'
$ws.Range('G13').Value = 75
$ws.Range('H13').Value = '◆開発'

# Row 14
$ws.Range('A14').Value = '2026-02-16 18:39:57'
$ws.Range('B14').Value = '【急募】ダウンロードスクリプト開発とBOXアップロード依頼'
$ws.Range('C14').Value = 'システム開発'
$ws.Range('D14').Value = '100,000 円 ~ 200,000 円 / 募集期間 1 日、取引期間 0 日'
$ws.Range('E14').Value = '期限情報なし'
$ws.Range('F14').Value = 'https://www.lancers.jp/work/detail/5492631'
$ws.Range('G14').Value = 68
$ws.Range('H14').Value = '◆開発'

# Row 15
$ws.Range('A15').Value = '2026-02-16 18:39:57'
$ws.Range('B15').Value = 'オンラインくじサイトの作成 簡易版可'
$ws.Range('C15').Value = 'システム開発'
$ws.Range('D15').Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range('E15').Value = '期限情報なし'
$ws.Range('F15').Value = 'https://www.lancers.jp/work/detail/5492891'
$ws.Range('G15').Value = 33
$ws.Range('H15').Value = '◇サイト'

# Row 16
$ws.Range('A16').Value = '2026-02-16 18:39:57'
$ws.Range('B16').Value = '【急募】ホテル公式LINEリニューアル運用サポート依頼'
$ws.Range('C16').Value = 'システム開発'
$ws.Range('D16').Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range('E16').Value = '期限情報なし'
$ws.Range('F16').Value = 'https://www.lancers.jp/work/detail/5492894'
$ws.Range('G16').Value = 18

# Row 17
$ws.Range('A17').Value = '2026-02-16 18:39:57'
$ws.Range('B17').Value = '【急募】生産計画の自動調整マクロ作成依頼'
$ws.Range('C17').Value = 'システム開発'
$ws.Range('D17').Value = '1,000 ~ 5,000 円 / 固定'
$ws.Range('E17').Value = '期限情報なし'
$ws.Range('F17').Value = 'https://www.lancers.jp/work/detail/5492925'
$ws.Range('G17').Value = 10

# Re-create hyperlinks for the URL column (F) and reapply the Hyperlink style
# so every F-cell keeps the same shared cell format as before.
$ws.Hyperlinks.Add($ws.Range('F2'), 'https://www.lancers.jp/work/detail/5460562')
$ws.Hyperlinks.Add($ws.Range('F3'), 'https://www.lancers.jp/work/detail/5473648')
$ws.Hyperlinks.Add($ws.Range('F4'), 'https://www.lancers.jp/work/detail/5492832')
$ws.Hyperlinks.Add($ws.Range('F5'), 'https://www.lancers.jp/work/detail/5217096')
$ws.Hyperlinks.Add($ws.Range('F6'), 'https://www.lancers.jp/work/detail/5460563')
$ws.Hyperlinks.Add($ws.Range('F7'), 'https://www.lancers.jp/work/detail/5488168')
$ws.Hyperlinks.Add($ws.Range('F8'), 'https://www.lancers.jp/work/detail/5492887')
$ws.Hyperlinks.Add($ws.Range('F9'), 'https://www.lancers.jp/work/detail/5492576')
$ws.Hyperlinks.Add($ws.Range('F10'), 'https://www.lancers.jp/work/detail/5492959')
$ws.Hyperlinks.Add($ws.Range('F11'), 'https://www.lancers.jp/work/detail/5492441')
$ws.Hyperlinks.Add($ws.Range('F12'), 'https://www.lancers.jp/work/detail/5492383')
$ws.Hyperlinks.Add($ws.Range('F13'), 'https://www.lancers.jp/work/detail/5468432')
$ws.Hyperlinks.Add($ws.Range('F14'), 'https://www.lancers.jp/work/detail/5492631')
$ws.Hyperlinks.Add($ws.Range('F15'), 'https://www.lancers.jp/work/detail/5492891')
$ws.Hyperlinks.Add($ws.Range('F16'), 'https://www.lancers.jp/work/detail/5492894')
$ws.Hyperlinks.Add($ws.Range('F17'), 'https://www.lancers.jp/work/detail/5492925')
$ws.Range('F2:F17').Style = 'Hyperlink'

